$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 2.75
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 1.62
